# Weekly update: insert two new price records (row 393-394) for
# "Clementina" (Especial / Primera) sourced from Región de O'Higgins,
# pushing the previously existing rows (393 onward) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 393 (Excel shifts the
# existing rows 393.. down to 395.. and copies formatting, including the
# date style in column D, from the row above).
$ws.Rows.Item(393).Insert()
$ws.Rows.Item(393).Insert()

# --- New row 393 ---
$ws.Cells.Item(393, 1).Value = 9
$ws.Cells.Item(393, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(393, 3).Value = "Metropolitana"
$ws.Cells.Item(393, 4).Value = 44516
$ws.Cells.Item(393, 5).Value = 13
$ws.Cells.Item(393, 6).Value = "Fruta"
$ws.Cells.Item(393, 7).Value = 100102
$ws.Cells.Item(393, 8).Value = "Cítricos"
$ws.Cells.Item(393, 9).Value = 100102004
$ws.Cells.Item(393, 10).Value = "Mandarina"
$ws.Cells.Item(393, 11).Value = "Clementina"
$ws.Cells.Item(393, 12).Value = "Especial"
$ws.Cells.Item(393, 13).Value = 400
$ws.Cells.Item(393, 14).Value = 11500
$ws.Cells.Item(393, 15).Value = 11500
$ws.Cells.Item(393, 16).Value = 11500
$ws.Cells.Item(393, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(393, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(393, 19).Value = 639
$ws.Cells.Item(393, 20).Value = 18

# --- New row 394 ---
$ws.Cells.Item(394, 1).Value = 9
$ws.Cells.Item(394, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(394, 3).Value = "Metropolitana"
$ws.Cells.Item(394, 4).Value = 44516
$ws.Cells.Item(394, 5).Value = 13
$ws.Cells.Item(394, 6).Value = "Fruta"
$ws.Cells.Item(394, 7).Value = 100102
$ws.Cells.Item(394, 8).Value = "Cítricos"
$ws.Cells.Item(394, 9).Value = 100102004
$ws.Cells.Item(394, 10).Value = "Mandarina"
$ws.Cells.Item(394, 11).Value = "Clementina"
$ws.Cells.Item(394, 12).Value = "Primera"
$ws.Cells.Item(394, 13).Value = 350
$ws.Cells.Item(394, 14).Value = 10000
$ws.Cells.Item(394, 15).Value = 10000
$ws.Cells.Item(394, 16).Value = 10000
$ws.Cells.Item(394, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(394, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(394, 19).Value = 556
$ws.Cells.Item(394, 20).Value = 18
